# General.xlsx — Crowdin translation sync
#
# 1) Two strings that contain Japanese furigana/ruby annotations (<rPh>) get
#    re-written without the reading-annotation text baked into the cell
#    value (the round-trip through the translation tool drops <phoneticPr>
#    entirely, so we restore the clean value here).
# 2) Two new rows are appended to the "Main" sheet with a new Male/Female
#    translation pair (EN label in column A, JA text in column B).
# 3) Rows 8 and 10 get their row height pinned explicitly (same values they
#    already rendered at) so the custom-height flag is persisted for them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clean up the CalloutCode4 / ErrorTooFar Japanese strings -------------
$ws.Range("B8").Value  = "~g~事案終了~s~`n現時点をもってパトロールに復帰せよ。"
$ws.Range("B10").Value = "遠すぎます。`nもっと近づいてください。"

# --- Pin the explicit row heights those two rows already use --------------
$ws.Rows(8).RowHeight  = 40.5
$ws.Rows(10).RowHeight = 27

# --- New Male/Female translation rows (rows 12 & 13) -----------------------
$ws.Range("A12").Value = "Male"
$ws.Range("B12").Value = "男性"
$ws.Range("A13").Value = "Female"
$ws.Range("B13").Value = "女性"
